$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# The paragraph right after "...denotes your name/s." currently holds two
# runs (" x " + "   ") that together read " x    ". Collapse them into a
# single run with that exact text (same run formatting as before).
$r1 = $d.Content
$r1.Find.Execute("denotes your name/s.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$p1 = $r1.Paragraphs(1)
$target1 = $p1.Next()
$tr1 = $target1.Range
$tr1.MoveEnd(1, -1) | Out-Null
# Force an actual content change (so the run-merge happens) even though
# the final text happens to equal the current combined text.
$tr1.Text = "__TMP__"
$tr1b = $target1.Range
$tr1b.MoveEnd(1, -1) | Out-Null
$tr1b.Text = " x    "

# --- Edit 2 -----------------------------------------------------------
# The paragraph right after "Produce a test suite for Requirements #6 and
# #7." currently holds a single run of seven spaces. Replace the text
# with "x" and tag the run with Spanish (es-ES) language.
$r2 = $d.Content
$r2.Find.Execute("Produce a test suite for Requirements #6 and #7.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$p2 = $r2.Paragraphs(1)
$target2 = $p2.Next()
$tr2 = $target2.Range
$tr2.MoveEnd(1, -1) | Out-Null
$tr2.Text = "x"
$tr2.LanguageID = "es-ES"
$tr2.LanguageIDFarEast = "es-ES"
$tr2.LanguageIDOther = "es-ES"
